{"js": "// Replace each three-digit x one-digit multiplication equation with its\n// updated version, as described in the commit diff. Each equation text is\n// unique within the document, so we can safely search for the exact old\n// string and replace it with the new one.\nconst replacements = [\n  [\"790\u00d72=1580\", \"814\u00d79=7326\"],\n  [\"997\u00d72=1994\", \"259\u00d76=1554\"],\n  [\"120\u00d75=600\", \"407\u00d78=3256\"],\n  [\"728\u00d73=2184\", \"541\u00d75=2705\"],\n  [\"699\u00d74=2796\", \"974\u00d79=8766\"],\n  [\"430\u00d76=2580\", \"390\u00d77=2730\"],\n  [\"738\u00d74=2952\", \"510\u00d76=3060\"],\n  [\"532\u00d72=1064\", \"736\u00d72=1472\"],\n  [\"711\u00d75=3555\", \"441\u00d77=3087\"],\n  [\"659\u00d76=3954\", \"205\u00d78=1640\"],\n  [\"343\u00d78=2744\", \"919\u00d73=2757\"],\n  [\"728\u00d76=4368\", \"843\u00d76=5058\"],\n  [\"405\u00d74=1620\", \"334\u00d79=3006\"],\n  [\"641\u00d75=3205\", \"946\u00d77=6622\"],\n  [\"611\u00d73=1833\", \"763\u00d76=4578\"],\n  [\"927\u00d76=5562\", \"864\u00d75=4320\"],\n  [\"382\u00d72=764\", \"674\u00d75=3370\"],\n  [\"350\u00d73=1050\", \"388\u00d73=1164\"],\n  [\"992\u00d77=6944\", \"326\u00d72=652\"],\n  [\"798\u00d75=3990\", \"498\u00d75=2490\"],\n  [\"147\u00d76=882\", \"125\u00d72=250\"],\n  [\"243\u00d76=1458\", \"201\u00d72=402\"],\n  [\"664\u00d79=5976\", \"231\u00d76=1386\"],\n  [\"163\u00d77=1141\", \"252\u00d77=1764\"],\n  [\"709\u00d73=2127\", \"856\u00d75=4280\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update each three-digit x one-digit multiplication equation to the\n# values produced by the latest generation run, per the commit diff.\n# Every equation string occurs exactly once in the document, so a simple\n# MatchCase / MatchWholeWord-off Find & ReplaceAll for each exact pair is\n# safe and deterministic.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"790\u00d72=1580\"; New = \"814\u00d79=7326\" },\n    @{ Old = \"997\u00d72=1994\"; New = \"259\u00d76=1554\" },\n    @{ Old = \"120\u00d75=600\"; New = \"407\u00d78=3256\" },\n    @{ Old = \"728\u00d73=2184\"; New = \"541\u00d75=2705\" },\n    @{ Old = \"699\u00d74=2796\"; New = \"974\u00d79=8766\" },\n    @{ Old = \"430\u00d76=2580\"; New = \"390\u00d77=2730\" },\n    @{ Old = \"738\u00d74=2952\"; New = \"510\u00d76=3060\" },\n    @{ Old = \"532\u00d72=1064\"; New = \"736\u00d72=1472\" },\n    @{ Old = \"711\u00d75=3555\"; New = \"441\u00d77=3087\" },\n    @{ Old = \"659\u00d76=3954\"; New = \"205\u00d78=1640\" },\n    @{ Old = \"343\u00d78=2744\"; New = \"919\u00d73=2757\" },\n    @{ Old = \"728\u00d76=4368\"; New = \"843\u00d76=5058\" },\n    @{ Old = \"405\u00d74=1620\"; New = \"334\u00d79=3006\" },\n    @{ Old = \"641\u00d75=3205\"; New = \"946\u00d77=6622\" },\n    @{ Old = \"611\u00d73=1833\"; New = \"763\u00d76=4578\" },\n    @{ Old = \"927\u00d76=5562\"; New = \"864\u00d75=4320\" },\n    @{ Old = \"382\u00d72=764\"; New = \"674\u00d75=3370\" },\n    @{ Old = \"350\u00d73=1050\"; New = \"388\u00d73=1164\" },\n    @{ Old = \"992\u00d77=6944\"; New = \"326\u00d72=652\" },\n    @{ Old = \"798\u00d75=3990\"; New = \"498\u00d75=2490\" },\n    @{ Old = \"147\u00d76=882\"; New = \"125\u00d72=250\" },\n    @{ Old = \"243\u00d76=1458\"; New = \"201\u00d72=402\" },\n    @{ Old = \"664\u00d79=5976\"; New = \"231\u00d76=1386\" },\n    @{ Old = \"163\u00d77=1141\"; New = \"252\u00d77=1764\" },\n    @{ Old = \"709\u00d73=2127\"; New = \"856\u00d75=4280\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    $found = $find.Execute([ref]$pair.Old, $false, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n\n    if (-not $found) {\n        Write-Host \"Could not find text to replace: $($pair.Old)\"\n    }\n}\n"}
